$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Config")

# Add new backup-folder config rows
$ws.Range("A15").Value = "vthhBackupFolder"
$ws.Range("B15").Value = "BackupFileOutput\VTHH"

$ws.Range("A16").Value = "nccBackupFolder"
$ws.Range("B16").Value = "BackupFileOutput\NCC"

$ws.Range("A17").Value = "mhBackupFolder"
$ws.Range("B17").Value = "BackupFileOutput\MuaHang"

# Match formatting of the existing KEY column cells (text number format)
$ws.Range("A15:A17").NumberFormat = "@"

# Re-affirm the "stage" value (row 7, column B) - condition for opening vsign file
$ws.Range("B7").Value = "1"

# Update the active cell/selection
$ws.Range("E12").Select()
